# auto commit by win-upload.bat
#
# Adds a new, blank "Title and Content" slide as slide 2 (right after the
# existing slide 1), mirroring a manual "New Slide" insert in the
# PowerPoint UI (Home > New Slide > Title and Content). The slide is left
# with empty placeholders, matching the freshly-inserted/unedited slide.

$p = $ppt.ActivePresentation

# ppLayoutText (2) == the "Title and Content" custom layout (slideLayout2.xml)
# used by the target slide: a title placeholder plus a single body/content
# placeholder (idx=1), both left blank.
$s = $p.Slides.Add(2, 2)
